# Refresh Kraken_Profits market-price snapshot values across all job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 1229.6
$ws.Range("J58").Value = 1899.3334
$ws.Range("L58").Value = 5698.0002
$ws.Range("N58").Value = -5998.0002
# Row 80
$ws.Range("H80").Value = 600
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 600
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1800
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3796
# Row 83
$ws.Range("H83").Value = 600
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 5400
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -15384
# Row 141
$ws.Range("H141").Value = 824
$ws.Range("I141").Value = 898
$ws.Range("J141").Value = 750
$ws.Range("K141").Value = 2694
$ws.Range("L141").Value = 2250
$ws.Range("M141").Value = 2486
$ws.Range("N141").Value = -12610

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 5033
$ws.Range("I2").Value = 2550
$ws.Range("K2").Value = 2550
$ws.Range("M2").Value = -2437
# Row 22
$ws.Range("H22").Value = 1508
$ws.Range("I22").Value = 516
$ws.Range("K22").Value = 516
$ws.Range("M22").Value = -217
# Row 45
$ws.Range("H45").Value = 1972.5
$ws.Range("I45").Value = 2080
$ws.Range("K45").Value = 2080
$ws.Range("M45").Value = -1703
# Row 61
$ws.Range("H61").Value = 5241
$ws.Range("I61").Value = 4754.857
$ws.Range("K61").Value = 4754.857
$ws.Range("M61").Value = -4542.857
# Row 74
$ws.Range("H74").Value = 4436.5713
$ws.Range("I74").Value = 2892.2354
$ws.Range("J74").Value = 11000
$ws.Range("K74").Value = 2892.2354
$ws.Range("L74").Value = 11000
$ws.Range("M74").Value = -2018.2354
$ws.Range("N74").Value = -12748
# Row 77
$ws.Range("H77").Value = 4436.5713
$ws.Range("I77").Value = 2892.2354
$ws.Range("J77").Value = 11000
$ws.Range("K77").Value = 14461.177
$ws.Range("L77").Value = 55000
$ws.Range("M77").Value = -10093.177
$ws.Range("N77").Value = -63736
# Row 116
$ws.Range("H116").Value = 5033
$ws.Range("I116").Value = 2550
$ws.Range("K116").Value = 2550
$ws.Range("M116").Value = -256
# Row 122
$ws.Range("H122").Value = 4231.933
$ws.Range("I122").Value = 4664.9165
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 13994.7495
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -11544.7495
$ws.Range("N122").Value = -12400
# Row 136
$ws.Range("H136").Value = 5241
$ws.Range("I136").Value = 4754.857
$ws.Range("K136").Value = 14264.571
$ws.Range("M136").Value = -11714.571

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 5033
$ws.Range("I3").Value = 2550
$ws.Range("K3").Value = 2550
$ws.Range("M3").Value = -2436
# Row 5
$ws.Range("H5").Value = 1098
$ws.Range("I5").Value = 897.5
$ws.Range("J5").Value = 1499
$ws.Range("K5").Value = 897.5
$ws.Range("L5").Value = 1499
$ws.Range("M5").Value = -784.5
$ws.Range("N5").Value = -1725
# Row 105
$ws.Range("H105").Value = 4448.5
$ws.Range("I105").Value = 3897.5
$ws.Range("J105").Value = 4999.5
$ws.Range("K105").Value = 3897.5
$ws.Range("L105").Value = 4999.5
$ws.Range("M105").Value = -2150.5
$ws.Range("N105").Value = -8493.5
# Row 107
$ws.Range("H107").Value = 2175
$ws.Range("J107").Value = 750
$ws.Range("L107").Value = 750
$ws.Range("N107").Value = -4590

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1772
$ws.Range("I58").Value = 912
$ws.Range("K58").Value = 912
$ws.Range("M58").Value = -709
# Row 136
$ws.Range("H136").Value = 1772
$ws.Range("I136").Value = 912
$ws.Range("K136").Value = 2736
$ws.Range("M136").Value = -186

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2500
$ws.Range("I5").Value = 1250
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 3750
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -3638
$ws.Range("N5").Value = -15224
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 131
$ws.Range("H131").Value = 1741.8572
$ws.Range("J131").Value = 2190.1428
$ws.Range("L131").Value = 6570.428400000001
$ws.Range("N131").Value = -16650.4284
# Row 135
$ws.Range("H135").Value = 2500
$ws.Range("I135").Value = 1250
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 11250
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -8715
$ws.Range("N135").Value = -50070
# Row 138
$ws.Range("H138").Value = 3588.2354
$ws.Range("I138").Value = 3588.2354
$ws.Range("K138").Value = 10764.7062
$ws.Range("M138").Value = -5624.706200000001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1352.1111
$ws.Range("J97").Value = 2381
$ws.Range("L97").Value = 2381
$ws.Range("N97").Value = -3373

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# Row 7
$ws.Range("H7").Value = 3392.1428
$ws.Range("I7").Value = 3049
$ws.Range("J7").Value = 4250
$ws.Range("K7").Value = 3049
$ws.Range("L7").Value = 4250
$ws.Range("M7").Value = -2937
$ws.Range("N7").Value = -4474
# Row 68
$ws.Range("H68").Value = 2138.4614
$ws.Range("I68").Value = 2088.889
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 2088.889
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -1339.889
$ws.Range("N68").Value = -3748
# Row 71
$ws.Range("H71").Value = 2138.4614
$ws.Range("I71").Value = 2088.889
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 10444.445
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -6700.445
$ws.Range("N71").Value = -18738
# Row 82
$ws.Range("H82").Value = 1602.6
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 13
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 13
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -735
# Row 85
$ws.Range("H85").Value = 1602.6
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 13
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 13
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -2509
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 122
$ws.Range("H122").Value = 3424.25
$ws.Range("I122").Value = 3870.2856
$ws.Range("J122").Value = 2799.8
$ws.Range("K122").Value = 11610.8568
$ws.Range("L122").Value = 8399.400000000001
$ws.Range("M122").Value = -9160.856800000001
$ws.Range("N122").Value = -13299.4
# Row 126
$ws.Range("H126").Value = 3392.1428
$ws.Range("I126").Value = 3049
$ws.Range("J126").Value = 4250
$ws.Range("K126").Value = 9147
$ws.Range("L126").Value = 12750
$ws.Range("M126").Value = -6677
$ws.Range("N126").Value = -17690

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 37353.273
$ws.Range("I2").Value = 40088.6
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 40088.6
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -39976.6
$ws.Range("N2").Value = -10224
# Row 132
$ws.Range("H132").Value = 6795.0713
$ws.Range("I132").Value = 5566.4546
$ws.Range("J132").Value = 11300
$ws.Range("K132").Value = 16699.3638
$ws.Range("L132").Value = 33900
$ws.Range("M132").Value = -14169.3638
$ws.Range("N132").Value = -38960
